$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Set G5 to the same value as F5/E5 ("OUI")
$ws.Range("G5").Value = $ws.Range("F5").Value2

# Update the active selection to G6, matching the saved view state
$ws.Activate()
$ws.Range("G6").Select()
